$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.904.05'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '1.781.15'
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.24'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("E6").Value = '  -1.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.79'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0679'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0936'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").Value = '2.037.61'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").Value = '1.871.68'
$ws.Range("E13").Value = '  +3.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.23'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.88%  '
$ws.Range("D15").Value = '33.907.80'
$ws.Range("E15").Value = '  -0.83%  '
$ws.Range("E16").Value = '  -3.51%  '
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.01'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.52'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.21%  '
$ws.Range("E20").Value = '  -2.46%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.61'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.03'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.97'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.10'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.90%  '
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("E28").Value = '  -1.00%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  +1.04%  '
$ws.Range("E31").Value = '  -3.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.60'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.79%  '
$ws.Range("E33").Value = '  -0.68%  '
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("D35").Value = '1.391.50'
$ws.Range("E35").Value = '  -1.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.638'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.46%  '
$ws.Range("E37").Value = '  -1.49%  '
$ws.Range("E38").Value = '  -1.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.27'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.94%  '
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("E41").Value = '  -2.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '78.42'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.57'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +12.63%  '
$ws.Range("E44").Value = '  -3.10%  '
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("E46").Value = '  +2.61%  '
$ws.Range("D47").Value = '0.0₆0135'
$ws.Range("E47").Value = '  +5.99%  '
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("D49").Value = '1.938.99'
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '105.62'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.56%  '
$ws.Range("E51").Value = '  -0.08%  '
